$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Target change (per the diff):
#   - The paragraph "Why amazon? Why compose? Describe in terms of NFRs."
#     is rewritten and split into two new paragraphs discussing why AWS and
#     Compose/MongoDB were chosen, while keeping the (hidden) "_GoBack"
#     bookmark positioned at the very end of the new second paragraph.
# ---------------------------------------------------------------------------

$para1Text = "Only a limited number of cloud providers were considered in the trade-off analysis. For security purposes, Amazon Web Services was selected as the cloud provider to host CometBites application. Its highly scalability allows the engineers to easily upgrade the servers hardware. Security groups are easy to configure and it is possible to filter ports to the web server which hosts the backend. In addition, the ELB (Elastic Load Balancer) providers an additional layer of protection to the application. It is also important to mention that the team is already experienced in using these services, another reason to select this provider since the project has strict time limitations. "

$para2Text = "The web servers hosted in the virtual machines are connected a MongoDB database. MongoDB is a good database option because it provides the functionality of storing the data in document-based collections, which allow different values for a same class. For instance, the Card class inside the database might have 4 attributes (name, number, cvv, and expiration date) if it is a credit card, but only 2 (name, and number) if it is a comet card. Finally, we chose to use Compose because it exposes MongoDB as a service and handles the operations internally, using a GUI to provide these features."

# Step 1: delete the old closing sentence "Describe in terms of NFRs." which
# sits right after the (hidden) _GoBack bookmark.
$findDescribe = $d.Content
$findDescribe.Find.Execute("Describe in terms of NFRs.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# Step 2: insert the new second-paragraph text immediately at the bookmark's
# (now collapsed) position using InsertAfter -- this places the new text
# *before* the bookmark, which keeps the bookmark anchored at the end of the
# paragraph (matching the target layout).
$bm = $d.Bookmarks("_GoBack")
$insPara2 = $d.Range($bm.Start, $bm.Start)
$insPara2.InsertAfter($para2Text)

# Step 3: replace the old opening sentence "Why amazon? Why compose? " with
# the new first-paragraph text.
$findWhy = $d.Content
$findWhy.Find.Execute("Why amazon? Why compose? ", $true, $false, $false, $false, $false, $true, 1, $false, $para1Text, 2) | Out-Null

# Step 4: split the single paragraph into two paragraphs by inserting a
# paragraph break right after the new first-paragraph text (and before the
# bookmark + new second-paragraph text).
$findPara1 = $d.Content
$findPara1.Find.Execute($para1Text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $d.Range($findPara1.End, $findPara1.End)
$splitPoint.InsertParagraphAfter()
